$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "iProc_TC_ID_1"
$ws.Cells.Item(2, 2).Value = "@Smoke Verify Elumina Login"
$ws.Cells.Item(2, 3).Value = "passed"

$ws.Cells.Item(3, 1).Value = "iProc_TC_ID_1A"
$ws.Cells.Item(3, 2).Value = "@Smoke Verify Elumina Registration"
$ws.Cells.Item(3, 3).Value = "passed"

$ws.Cells.Item(4, 1).Value = "iProc_TC_ID_3"
$ws.Cells.Item(4, 2).Value = "@Smoke Verify CandidatesInvExam"
$ws.Cells.Item(4, 3).Value = "passed"

$ws.Cells.Item(5, 1).Value = "iProc_TC_ID_5"
$ws.Cells.Item(5, 2).Value = "@SmokeValidation of user authentication by valid Candidate Creadentials`n"
$ws.Cells.Item(5, 3).Value = "passed"
$ws.Rows.Item(5).AutoFit()

$ws.Cells.Item(6, 1).Value = "iProc_TC_ID_8"
$ws.Cells.Item(6, 2).Value = "@Smoke Validation of candidate choosing proctoring exam in dashboard"
$ws.Cells.Item(6, 3).Value = "passed"

$ws.Cells.Item(7, 1).Value = "iProc_TC_ID_23"
$ws.Cells.Item(7, 2).Value = "@iProctorlink Verify CandidatesExam"
$ws.Cells.Item(7, 3).Value = "passed"

$ws.Cells.Item(8, 1).Value = "iProc_TC_ID_25"
$ws.Cells.Item(8, 2).Value = "@iProctorlink Verify CandidatesExam"
$ws.Cells.Item(8, 3).Value = "passed"

$ws.Cells.Item(9, 1).Value = "iProc_TC_ID_28"
$ws.Cells.Item(9, 2).Value = "@iProctorlink Verify Elumina Invigilator Dashboard"
$ws.Cells.Item(9, 3).Value = "passed"

$ws.Cells.Item(10, 1).Value = "iProc_TC_ID_56"
$ws.Cells.Item(10, 2).Value = "@Smoke Verify Validation of Invigilator Dashboard Proctor"
$ws.Cells.Item(10, 3).Value = "passed"

$ws.Cells.Item(11, 1).Value = "iProc_TC_ID_57"
$ws.Cells.Item(11, 2).Value = "@Smoke Validation of Navigating to an exam from the dashboard to invigilate"
$ws.Cells.Item(11, 3).Value = "passed"

$ws.Cells.Item(12, 1).Value = "iProc_TC_ID_58"
$ws.Cells.Item(12, 2).Value = "@iProctorlink Verify Validation of `"Start Exam`" (All Candidates) Proctor "
$ws.Cells.Item(12, 3).Value = "timedOut"

$ws.Cells.Item(13, 1).Value = "iProc_TC_ID_59"
$ws.Cells.Item(13, 2).Value = "@iProctorlink Verify Validation of `"Lock Exam`" from Live monitor Proctor"
$ws.Cells.Item(13, 3).Value = "passed"

$ws.Cells.Item(14, 1).Value = "iProc_TC_ID_61"
$ws.Cells.Item(14, 2).Value = "@iProctorlink Verify Validation of `"Resume Exam`" from Live monitor Proctor"
$ws.Cells.Item(14, 3).Value = "passed"

$ws.Cells.Item(15, 1).Value = "iProc_TC_ID_64"
$ws.Cells.Item(15, 2).Value = "@Smoke Verify Validation of `"Mark Attendance`" (All Candidates) Proctor"
$ws.Cells.Item(15, 3).Value = "failed"

$ws.Cells.Item(16, 1).Value = "iProc_TC_ID_70"
$ws.Cells.Item(16, 2).Value = "@Smoke Validation of Questions answered / Inprogress on the RHS of the Candidate page"
$ws.Cells.Item(16, 3).Value = "passed"

$ws.Cells.Item(17, 1).Value = "iProc_TC_ID_71"
$ws.Cells.Item(17, 2).Value = "@Smoke Validation of all the events generated on the RHS of the Candidate page"
$ws.Cells.Item(17, 3).Value = "passed"

$ws.Cells.Item(18, 1).Value = "Exam_Prerequisit_ID_01"
$ws.Cells.Item(18, 2).Value = "@iProctorlink Verify Create Exam With Content Section and Content Section Page"
$ws.Cells.Item(18, 3).Value = "failed"

$ws.Cells.Item(19, 1).Value = "Reg_Prerequisit_ID_01A"
$ws.Cells.Item(19, 2).Value = "@iProctorlink Verify Elumina Registration"
$ws.Cells.Item(19, 3).Value = "failed"

$ws.Cells.Item(20, 1).Value = "iProc_TC_ID_40"
$ws.Cells.Item(20, 2).Value = "@iProctorlink Verify Elumina Invigilator Dashboard"
$ws.Cells.Item(20, 3).Value = "timedOut"

$ws.Cells.Item(21, 1).Value = "@Smoke Verify Elumina Login and Create Exam"
$ws.Cells.Item(21, 3).Value = "passed"

$ws.Cells.Item(22, 1).Value = "@Smoke Verify Elumina RegistrationInv and add User and Invigilator"
$ws.Cells.Item(22, 3).Value = "passed"

$ws.Cells.Item(23, 1).Value = "iProc_TC_ID_11"
$ws.Cells.Item(23, 2).Value = "@Smoke Validation of `"I Authorise`" checkbox - To access Webcam, Microphone & Terms & Condition"
$ws.Cells.Item(23, 3).Value = "passed"

$ws.Cells.Item(24, 1).Value = "iProc_TC_ID_34"
$ws.Cells.Item(24, 2).Value = "@Smoke Validation of Exam section page (Offline Exam validation)"
$ws.Cells.Item(24, 3).Value = "passed"

$ws.Cells.Item(25, 1).Value = "iProc_TC_ID_39"
$ws.Cells.Item(25, 2).Value = "@Smoke Validation of submitting when the Candidate has not answered all Questions"
$ws.Cells.Item(25, 3).Value = "passed"
